# Actualización automática 2025-06-09 09:40:09
#
# Registers new "junio" (June) sales for two clients of the asesor
# "CASTRO ALCIVAR EDA MARIA" and ripples the totals through the three
# sheets of the workbook:
#   1. VENTAS POR GRUPO      - per-client sales broken out by product group
#   2. VENTA MENSUAL         - per-client sales broken out by month
#   3. CUMPLIMIENTO MENSUAL  - per-group budget vs. sales roll-up

$wb = $excel.ActiveWorkbook

$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")
$wsCumplimiento = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# --- 1. VENTAS POR GRUPO --------------------------------------------------
# SALAZAR BALLADARES MARIA ANGELICA sold 1037.8 of "240X80 PORCELANATO" (col D)
$wsGrupo.Range("D44").Value = 1037.8
# TAPIA BALLADARES LORENA ELIZABETH sold 1459.81 of "PORCELANATO" (col L)
$wsGrupo.Range("L51").Value = 1459.81

# Row 56 is the per-group "clientes cumplidos" counter ("X de 54"); both
# groups just gained one more client that reached goal this month.
$wsGrupo.Range("D56").Value = "2 de 54"
$wsGrupo.Range("L56").Value = "2 de 54"

# --- 2. VENTA MENSUAL ------------------------------------------------------
# Column F is "junio". Same two client sales recorded against that month.
$wsMensual.Range("F44").Value = 1037.8
$wsMensual.Range("F51").Value = 1459.81
# Row 56 totals roll up the whole column.
$wsMensual.Range("F56").Value = 16808.52

# --- 3. CUMPLIMIENTO MENSUAL -----------------------------------------------
# Column F ("CUMPLIMIENTO") width shrinks from 26 to 25 stored units.
# Excel's ColumnWidth property is offset from the stored OOXML width by
# 5/6 of a character, so subtract that to land on the exact target.
$wsCumplimiento.Columns.Item(6).ColumnWidth = 25 - (5/6)

# Row 3: "240X80 PORCELANATO" group picks up the 1037.8 sale.
$wsCumplimiento.Range("D3").Value = 5644.48
$wsCumplimiento.Range("E3").Value = 8083.52
$wsCumplimiento.Range("F3").Value = 0.4111655011655012

# Row 16: "PORCELANATO" group picks up the 1459.81 sale.
$wsCumplimiento.Range("D16").Value = 1806.28
$wsCumplimiento.Range("E16").Value = 43939.409
$wsCumplimiento.Range("F16").Value = 0.03948525073040216

# Row 19: TOTAL row reflects both new sales.
$wsCumplimiento.Range("D19").Value = 16904.28
$wsCumplimiento.Range("E19").Value = 74059.049
$wsCumplimiento.Range("F19").Value = 0.1858362065882615
